$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-pulled "dSF" (column F) values for rows 2-14
$values = @{
    2  = -2
    3  = 1
    4  = -2
    5  = 6
    6  = 3
    7  = 1
    8  = 4
    9  = 4
    10 = 4
    11 = -3
    12 = 3
    13 = 3
    14 = -4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
